$wb = $excel.ActiveWorkbook

# Fix underscore spacing in sheet names (remove stray space before the
# trailing "_Dataset"/"_Results" suffix).
$wb.Worksheets.Item(1).Name = "ASD>CTRL_DEGs_Dataset"
$wb.Worksheets.Item(2).Name = "ASD>CTRL_DEGs_Results"
$wb.Worksheets.Item(3).Name = "ASD<CTRL_DEGs_Dataset"
$wb.Worksheets.Item(4).Name = "ASD<CTRL_DEGs_Results"

# Update each sheet's selected/active cell, finishing on the sheet that
# should end up as the active tab (sheet 4).
$wb.Worksheets.Item(1).Range("H28").Select()
$wb.Worksheets.Item(2).Range("C15").Select()
$wb.Worksheets.Item(3).Range("K33").Select()
$wb.Worksheets.Item(4).Range("L34").Select()
